$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Swap the match data held in rows 40 and 41 (everything except
#    the row index in column A and the timestamp columns K/M/O/Q/S/U,
#    which stay put).
# ---------------------------------------------------------------
$cols = @("F","G","H","I","J","L","N","P","R","T","V")

$row40 = @{}
$row41 = @{}
foreach ($c in $cols) {
    $row40[$c] = $ws.Range($c + "40").Value2
    $row41[$c] = $ws.Range($c + "41").Value2
}

foreach ($c in $cols) {
    $ws.Range($c + "40").Value = $row41[$c]
    $ws.Range($c + "41").Value = $row40[$c]
}

# ---------------------------------------------------------------
# 2) Append a new row 43 for the G.A. Eagles vs Sittard match,
#    copying the formatting of an existing data row so the styles
#    (bold/bordered index column, date-formatted match-date column)
#    line up with the rest of the table.
# ---------------------------------------------------------------
$ws.Range("A2:V2").Copy($ws.Range("A43:V43"))

$ws.Range("A43").Value = 42
$ws.Range("B43").Value = "netherlands"
$ws.Range("C43").Value = "eredivisie"
$ws.Range("D43").Value = "2023-2024"
$ws.Range("E43").Value = 45191.83333333334
$ws.Range("F43").Value = "G.A. Eagles"
$ws.Range("G43").Value = 3
$ws.Range("H43").Value = "Sittard"
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 1.75
$ws.Range("K43").Value = "17/09/2023 11:43"
$ws.Range("L43").Value = 2.38
$ws.Range("M43").Value = "22/09/2023 19:54"
$ws.Range("N43").Value = 4.16
$ws.Range("O43").Value = "17/09/2023 11:43"
$ws.Range("P43").Value = 3.6
$ws.Range("Q43").Value = "22/09/2023 19:59"
$ws.Range("R43").Value = 4.33
$ws.Range("S43").Value = "17/09/2023 11:43"
$ws.Range("T43").Value = 3.03
$ws.Range("U43").Value = "22/09/2023 19:54"
$ws.Range("V43").Value = "https://www.betexplorer.com/football/netherlands/eredivisie/g-a-eagles-sittard/rDgC2rpF/"
